$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for "United Russia" party (rows 16-22)
$dates = @(42583, 43070, 43313, 43647, 43800, 43862, 44044)
$values = @(50, 48, 44, 44, 47, 45, 45)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item(2, 1).NumberFormat
    $ws.Cells.Item($row, 2).Value = "United Russia"
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

# Update the view: scroll position and selection (matches Excel's saved view state)
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("C10").Select()
